$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.165.41"
Set-TextValue $ws.Range("E2") "  -2.79%  "
Set-TextValue $ws.Range("D3") "3.170.94"
Set-TextValue $ws.Range("E3") "  -7.95%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "561.86"
Set-TextValue $ws.Range("E5") "  -4.07%  "
Set-TextValue $ws.Range("D6") "169.92"
Set-TextValue $ws.Range("E6") "  -2.46%  "
Set-TextValue $ws.Range("B7") "XRP"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D7") "0.611"
Set-TextValue $ws.Range("E7") "  +1.44%  "
Set-TextValue $ws.Range("B8") "USDC"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "3.167.37"
Set-TextValue $ws.Range("E9") "  -8.01%  "
Set-TextValue $ws.Range("E10") "  -6.33%  "
Set-TextValue $ws.Range("D11") "6.60"
Set-TextValue $ws.Range("E11") "  -5.08%  "
Set-TextValue $ws.Range("E12") "  -5.51%  "
Set-TextValue $ws.Range("D13") "3.720.45"
Set-TextValue $ws.Range("E13") "  -7.99%  "
Set-TextValue $ws.Range("E14") "  +1.09%  "
Set-TextValue $ws.Range("D15") "27.33"
Set-TextValue $ws.Range("E15") "  -6.35%  "
Set-TextValue $ws.Range("D16") "64.186.66"
Set-TextValue $ws.Range("E16") "  -2.72%  "
Set-TextValue $ws.Range("D17") "0.0000162"
Set-TextValue $ws.Range("E17") "  -5.63%  "
Set-TextValue $ws.Range("D18") "3.174.45"
Set-TextValue $ws.Range("E18") "  -7.84%  "
Set-TextValue $ws.Range("D19") "5.70"
Set-TextValue $ws.Range("E19") "  -4.44%  "
Set-TextValue $ws.Range("D20") "13.04"
Set-TextValue $ws.Range("E20") "  -5.82%  "
Set-TextValue $ws.Range("D21") "352.74"
Set-TextValue $ws.Range("E21") "  -4.90%  "
Set-TextValue $ws.Range("D22") "7.18"
Set-TextValue $ws.Range("E22") "  -5.62%  "
Set-TextValue $ws.Range("D23") "1.00"
Set-TextValue $ws.Range("E23") "  +0.04%  "
Set-TextValue $ws.Range("D24") "69.10"
Set-TextValue $ws.Range("E24") "  -4.74%  "
Set-TextValue $ws.Range("D25") "0.501"
Set-TextValue $ws.Range("E25") "  -6.07%  "
Set-TextValue $ws.Range("E26") "  -3.59%  "
Set-TextValue $ws.Range("D27") "9.56"
Set-TextValue $ws.Range("E27") "  -1.74%  "
Set-TextValue $ws.Range("E28") "  -1.96%  "
Set-TextValue $ws.Range("E29") "  +0.39%  "
Set-TextValue $ws.Range("D30") "5.61"
Set-TextValue $ws.Range("E30") "  -3.22%  "
Set-TextValue $ws.Range("E31") "  -0.21%  "
Set-TextValue $ws.Range("E32") "  -4.60%  "
Set-TextValue $ws.Range("D33") "22.09"
Set-TextValue $ws.Range("E33") "  -6.51%  "
Set-TextValue $ws.Range("D34") "6.64"
Set-TextValue $ws.Range("E34") "  -5.55%  "
Set-TextValue $ws.Range("E35") "  -5.76%  "
Set-TextValue $ws.Range("B36") "Monero"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D36") "156.13"
Set-TextValue $ws.Range("E36") "  -3.36%  "
Set-TextValue $ws.Range("B37") "ImmutableX"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "1.44"
Set-TextValue $ws.Range("E37") "  -6.37%  "
Set-TextValue $ws.Range("D38") "0.809"
Set-TextValue $ws.Range("E38") "  -7.99%  "
Set-TextValue $ws.Range("D39") "25.94"
Set-TextValue $ws.Range("E39") "  -8.48%  "
Set-TextValue $ws.Range("D40") "2.53"
Set-TextValue $ws.Range("E40") "  -2.92%  "
Set-TextValue $ws.Range("D41") "1.70"
Set-TextValue $ws.Range("E41") "  -4.65%  "
Set-TextValue $ws.Range("D42") "2.602.54"
Set-TextValue $ws.Range("E42") "  -6.79%  "
Set-TextValue $ws.Range("D43") "4.17"
Set-TextValue $ws.Range("E43") "  -6.68%  "
Set-TextValue $ws.Range("E44") "  -6.59%  "
Set-TextValue $ws.Range("B45") "Hedera"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D45") "0.0655"
Set-TextValue $ws.Range("E45") "  -4.77%  "
Set-TextValue $ws.Range("B46") "Bittensor"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D46") "329.38"
Set-TextValue $ws.Range("E46") "  +0.39%  "
Set-TextValue $ws.Range("B47") "OKB"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D47") "38.84"
Set-TextValue $ws.Range("E47") "  -2.46%  "
Set-TextValue $ws.Range("B48") "InjectiveProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "23.92"
Set-TextValue $ws.Range("E48") "  -5.18%  "
Set-TextValue $ws.Range("E49") "  -7.34%  "
Set-TextValue $ws.Range("D50") "0.102"
Set-TextValue $ws.Range("E50") "  -0.99%  "
Set-TextValue $ws.Range("E51") "  -0.04%  "
